$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$a2 = @'
https://p16-images-sign-sg.tokopedia-static.net/tos-alisg-i-aphluv4xwc-sg/68c57dc5515e456bb2e5de9363826388~tplv-aphluv4xwc-white-pad-v1:1600:1600.jpeg?lk3s=0ccea506&x-expires=1770020675&x-signature=1n%2F3G7DjLHMPtmH2bweKZ9v%2BYCw%3D&x-signature-webp=FH951L5gXfAQcpeXyNwcqBJlbO8%3D
'@

$b2 = @'
ASUS GAMING K16 K3605VC GeForce RTX 3050 - i5 13420H RAM 16GB 512GB SSD - RAM 8GB, TANPA ANTIGORES
'@

$c2 = @'
Rp10.348.899
'@

$d2 = @'
RTX 3000
The Ultimate Play GeForce RTX 30 Series with NVIDIA DLSS & Ray Tracing
ASUS VIVOBOOK 16X K3605VC -I535B6T-HM
I5 13420H RTX3050 4GB/ 8GB / 16GB 512GB W11+OHS+O365 16.0 WUXGA BLACK
Specification :
Graphics : NVIDIA GeForce RTX 3050 4GB Laptop GPU 4GB GDDR6, Intel UHD Graphics
Processor : Intel Core i5-13420H Processor 2.1 GHz (12MB Cache, up to 4.6 GHz, 8 cores, 12 Threads)
Operating System : Windows 11 Home
Microsoft Office :Office Home and Student 2021 included
Display : 16.0-inch, WUXGA (1920 x 1200) 16:10 aspect ratio, IPS-level Panel, LED Backlit, 144Hz refresh rate, 300nits, 45% NTSC color gamut, Anti-glare display, TÜV Rheinland-certified, (Screen-to-body ratio)86%
Memory :  8GB DDR4 on board + 8GB DDR4 SO-DIMM ( tersedia 2 varian 8GB & 16GB )
Max Total system memory up to:24GB
Storage : 512GB M.2 NVMe PCIe 4.0 SSD (replace for upgrade)
Ports :
1x USB 3.2 Gen 1 Type-C support power delivery
2x USB 3.2 Gen 1 Type-A
1x HDMI 2.1 TMDS
1x 3.5mm Combo Audio Jack
1x DC-in
SD 4.0 card reader
Camera : 720p HD camera With privacy shutter
Network and Communication : Wi-Fi 6E(802.11ax) (Triple band) 2*2 + Bluetooth 5.3 Wireless Card (*Bluetooth version may change with OS version different.)
Battery : 50WHrs, 3S1P, 3-cell Li-ion
Power Supply : 120W AC Adapter,
Weight : 1.80 kg
Keyboard & Touchpad : Backlit Chiclet Keyboard with Num-key
Free backpak asus
Warranty By ASUS Indonesia 2 Years
===========================================================
KETENTUAN TOKO
- Buka senin-sabtu (10.00-17.00) Minggu/tgl merah (chat admin)
- Garansi toko 5 hari wajib video unboxing dan belum aktivasi office (AKTIVASI OFFICE GARANSI TOKO HANGUS)
- Seluruh laptop di NVIDIA OFFICIAL merupakan GARANSI RESMI
- Faktur pajak, chat admin untuk s&k nya
'@

$ws.Range("A2").Value = $a2
$ws.Range("B2").Value = $b2
$ws.Range("C2").Value = $c2
$ws.Range("D2").Value = $d2
